$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimación de Procesos")

# New header / value pairs to add in columns F (Proceso) and G (Vista)
# keyed by row number.
$data = @{
    5  = @{ F = "Solicitud" }
    6  = @{ F = "N/A"; G = "N/A" }
    8  = @{ F = "N/A"; G = "N/A" }
    9  = @{ F = "Seguridad"; G = "CoordinadoresNomina" }
    10 = @{ F = "Solicitud"; G = "Captura" }
    11 = @{ F = "Solicitud"; G = "CargaDocumento" }
    12 = @{ F = "Solicitud"; G = "CargaMasiva" }
    14 = @{ F = "Nomina"; G = "Periodos" }
    15 = @{ F = "Solicitud"; G = "Administracion" }
    16 = @{ F = "N/A"; G = "N/A" }
    17 = @{ F = "Solicitud"; G = "Notificaciones" }
    18 = @{ F = "Proceso"; G = "GeneracionEXCEL" }
    21 = @{ F = "Reporte"; G = "General" }
    22 = @{ F = "Reporte"; G = "Auditoria" }
    23 = @{ F = "Proceso"; G = "CierreMovimientos" }
}

foreach ($row in $data.Keys) {
    $pair = $data[$row]
    if ($pair.ContainsKey("F")) {
        $ws.Cells.Item($row, 6).Value = $pair["F"]
    }
    if ($pair.ContainsKey("G")) {
        $ws.Cells.Item($row, 7).Value = $pair["G"]
    }
}

# Column width adjustments: split the old E:H uniform width block so that
# F ("Proceso") and G ("Vista") get bestFit custom widths (closest width
# achievable given the host's column-width quantization).
$ws.Columns.Item(6).ColumnWidth = 18.666666666666668
$ws.Columns.Item(7).ColumnWidth = 25.166666666666668

# Update the active selection on the sheet to match the saved state.
$ws.Range("G13").Select()
